$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel, by default, auto-converts a numeric/date-looking string typed
# into a cell's .Value into a real number/date. The source cells here are
# plain text (t="inlineStr") holding things like "2025/11/19" and "4.45",
# so we must force the new value to stay literal text. Temporarily
# flipping NumberFormat to "@" (Text) makes the assignment keep the
# literal string instead of parsing it; resetting the style back to
# "Normal" afterwards removes that temporary formatting stamp so the
# cell's style is left exactly as it was before the edit.
function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "A2" "2025/11/20"
Set-TextValue "B2" "4.47"

Set-TextValue "A8" "2025/11/20"
Set-TextValue "B8" "7.43"

Set-TextValue "A14" "2025/11/20"
Set-TextValue "B14" "2.78"

Set-TextValue "A20" "2025/11/20"
Set-TextValue "B20" "11.87"

Set-TextValue "A26" "2025/11/20"
Set-TextValue "B26" "9.61"

Set-TextValue "A32" "2025/11/20"
Set-TextValue "B32" "24.84"

Set-TextValue "A38" "2025/11/20"

Set-TextValue "A44" "2025/11/20"
Set-TextValue "B44" "9.90"

Set-TextValue "A50" "2025/11/20"
Set-TextValue "B50" "11.11"

Set-TextValue "A56" "2025/11/20"
Set-TextValue "B56" "32.78"

Set-TextValue "A62" "2025/11/20"
Set-TextValue "B62" "10.92"

Set-TextValue "A68" "2025/11/20"
Set-TextValue "B68" "12.21"

Set-TextValue "A74" "2025/11/20"
Set-TextValue "B74" "15.12"
